$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4766.7144
$ws.Range("J17").Value = 4766.7144
$ws.Range("L17").Value = 14300.1432
$ws.Range("N17").Value = -14636.1432
$ws.Range("H86").Value = 4184.0713
$ws.Range("J86").Value = 4727.8
$ws.Range("L86").Value = 4727.8
$ws.Range("N86").Value = -6973.8
$ws.Range("H89").Value = 4184.0713
$ws.Range("J89").Value = 4727.8
$ws.Range("L89").Value = 23639
$ws.Range("N89").Value = -34871
$ws.Range("H131").Value = 9336.467000000001
$ws.Range("I131").Value = 1079.2858
$ws.Range("K131").Value = 3237.8574
$ws.Range("M131").Value = 1802.1426
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5165.485
$ws.Range("I32").Value = 4885.839
$ws.Range("K32").Value = 4885.839
$ws.Range("M32").Value = -4598.839
$ws.Range("H61").Value = 4651.077
$ws.Range("I61").Value = 1781.4286
$ws.Range("K61").Value = 1781.4286
$ws.Range("M61").Value = -1569.4286
$ws.Range("H102").Value = 1244
$ws.Range("J102").Value = 994
$ws.Range("L102").Value = 994
$ws.Range("N102").Value = -4238
$ws.Range("H104").Value = 39998.5
$ws.Range("J104").Value = 39998.5
$ws.Range("L104").Value = 39998.5
$ws.Range("N104").Value = -46986.5
$ws.Range("H110").Value = 835.5
$ws.Range("J110").Value = 775
$ws.Range("L110").Value = 775
$ws.Range("N110").Value = -4865
$ws.Range("H132").Value = 1046.5
$ws.Range("I132").Value = 869.4783
$ws.Range("K132").Value = 2608.4349
$ws.Range("M132").Value = -78.4349000000002
$ws.Range("H136").Value = 4651.077
$ws.Range("I136").Value = 1781.4286
$ws.Range("K136").Value = 5344.2858
$ws.Range("M136").Value = -2794.2858
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1923.95
$ws.Range("I31").Value = 1707.6111
$ws.Range("J31").Value = 3871
$ws.Range("K31").Value = 1707.6111
$ws.Range("L31").Value = 3871
$ws.Range("M31").Value = -1412.6111
$ws.Range("N31").Value = -4461
$ws.Range("H34").Value = 1923.95
$ws.Range("I34").Value = 1707.6111
$ws.Range("J34").Value = 3871
$ws.Range("K34").Value = 1707.6111
$ws.Range("L34").Value = 3871
$ws.Range("M34").Value = -1505.6111
$ws.Range("N34").Value = -4275
$ws.Range("H58").Value = 3386.25
$ws.Range("I58").Value = 3348.5
$ws.Range("K58").Value = 3348.5
$ws.Range("M58").Value = -3145.5
$ws.Range("H86").Value = 3092.1428
$ws.Range("I86").Value = 2692.6667
$ws.Range("J86").Value = 3391.75
$ws.Range("K86").Value = 2692.6667
$ws.Range("L86").Value = 3391.75
$ws.Range("M86").Value = -1569.6667
$ws.Range("N86").Value = -5637.75
$ws.Range("H89").Value = 3092.1428
$ws.Range("I89").Value = 2692.6667
$ws.Range("J89").Value = 3391.75
$ws.Range("K89").Value = 13463.3335
$ws.Range("L89").Value = 16958.75
$ws.Range("M89").Value = -7847.333500000001
$ws.Range("N89").Value = -28190.75
$ws.Range("H132").Value = 1144.0769
$ws.Range("I132").Value = 1082.4
$ws.Range("J132").Value = 1349.6666
$ws.Range("K132").Value = 3247.2
$ws.Range("L132").Value = 4048.9998
$ws.Range("M132").Value = -717.2000000000003
$ws.Range("N132").Value = -9108.9998
$ws.Range("H136").Value = 3386.25
$ws.Range("I136").Value = 3348.5
$ws.Range("K136").Value = 10045.5
$ws.Range("M136").Value = -7495.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27500078
$ws.Range("I4").Value = 27500078
$ws.Range("K4").Value = 82500234
$ws.Range("M4").Value = -82500122
$ws.Range("H68").Value = 15636750
$ws.Range("I68").Value = 20999
$ws.Range("J68").Value = 20842000
$ws.Range("K68").Value = 62997
$ws.Range("L68").Value = 62526000
$ws.Range("M68").Value = -62186
$ws.Range("N68").Value = -62527622
$ws.Range("H71").Value = 15636750
$ws.Range("I71").Value = 20999
$ws.Range("J71").Value = 20842000
$ws.Range("K71").Value = 188991
$ws.Range("L71").Value = 187578000
$ws.Range("M71").Value = -184935
$ws.Range("N71").Value = -187586112
$ws.Range("H113").Value = 946.7727
$ws.Range("I113").Value = 1421.75
$ws.Range("J113").Value = 841.2222
$ws.Range("K113").Value = 4265.25
$ws.Range("L113").Value = 2523.6666
$ws.Range("M113").Value = -2095.25
$ws.Range("N113").Value = -6863.6666
$ws.Range("H121").Value = 4738.706
$ws.Range("I121").Value = 2866.2
$ws.Range("J121").Value = 5518.9165
$ws.Range("K121").Value = 8598.599999999999
$ws.Range("L121").Value = 16556.7495
$ws.Range("M121").Value = -7288.599999999999
$ws.Range("N121").Value = -19176.7495
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 3189.6667
$ws.Range("I19").Value = 2030
$ws.Range("J19").Value = 3769.5
$ws.Range("K19").Value = 2030
$ws.Range("L19").Value = 3769.5
$ws.Range("M19").Value = -1742
$ws.Range("N19").Value = -4345.5
$ws.Range("H80").Value = 3145.5715
$ws.Range("I80").Value = 2331
$ws.Range("K80").Value = 2331
$ws.Range("M80").Value = -1333
$ws.Range("H83").Value = 3145.5715
$ws.Range("I83").Value = 2331
$ws.Range("K83").Value = 11655
$ws.Range("M83").Value = -6663
$ws.Range("H113").Value = 1646.4
$ws.Range("I113").Value = 2069
$ws.Range("J113").Value = 1012.5
$ws.Range("K113").Value = 2069
$ws.Range("L113").Value = 1012.5
$ws.Range("M113").Value = 101
$ws.Range("N113").Value = -5352.5
$ws.Range("H132").Value = 1748.5294
$ws.Range("I132").Value = 2016.2142
$ws.Range("K132").Value = 6048.642599999999
$ws.Range("M132").Value = -3518.642599999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8009.524
$ws.Range("I7").Value = 4303.25
$ws.Range("J7").Value = 8881.588
$ws.Range("K7").Value = 4303.25
$ws.Range("L7").Value = 8881.588
$ws.Range("M7").Value = -4191.25
$ws.Range("N7").Value = -9105.588
$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1280
$ws.Range("H68").Value = 3716.8333
$ws.Range("J68").Value = 3716.8333
$ws.Range("L68").Value = 3716.8333
$ws.Range("N68").Value = -5214.8333
$ws.Range("H71").Value = 3716.8333
$ws.Range("J71").Value = 3716.8333
$ws.Range("L71").Value = 18584.1665
$ws.Range("N71").Value = -26072.1665
$ws.Range("H82").Value = 4099.5
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 5199
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 5199
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -5921
$ws.Range("H85").Value = 4099.5
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 5199
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 5199
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -7695
$ws.Range("H87").Value = 24000
$ws.Range("J87").Value = 24000
$ws.Range("L87").Value = 24000
$ws.Range("N87").Value = -26246
$ws.Range("H90").Value = 24000
$ws.Range("J90").Value = 24000
$ws.Range("L90").Value = 72000
$ws.Range("N90").Value = -83232
$ws.Range("H126").Value = 8009.524
$ws.Range("I126").Value = 4303.25
$ws.Range("J126").Value = 8881.588
$ws.Range("K126").Value = 12909.75
$ws.Range("L126").Value = 26644.764
$ws.Range("M126").Value = -10439.75
$ws.Range("N126").Value = -31584.764
$ws.Range("H136").Value = 3411.5557
$ws.Range("I136").Value = 1814.8572
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 5444.571599999999
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -2894.571599999999
$ws.Range("N136").Value = -32100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2000000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H62").Value = 5088.1113
$ws.Range("I62").Value = 3749.3333
$ws.Range("J62").Value = 5757.5
$ws.Range("K62").Value = 3749.3333
$ws.Range("L62").Value = 5757.5
$ws.Range("M62").Value = -3125.3333
$ws.Range("N62").Value = -7005.5
$ws.Range("H65").Value = 5088.1113
$ws.Range("I65").Value = 3749.3333
$ws.Range("J65").Value = 5757.5
$ws.Range("K65").Value = 18746.6665
$ws.Range("L65").Value = 28787.5
$ws.Range("M65").Value = -15626.6665
$ws.Range("N65").Value = -35027.5
$ws.Range("H122").Value = 938.8182
$ws.Range("I122").Value = 938.8182
$ws.Range("K122").Value = 2816.4546
$ws.Range("M122").Value = -366.4546
$ws.Range("H136").Value = 1945.3334
$ws.Range("I136").Value = 1666
$ws.Range("K136").Value = 4998
$ws.Range("M136").Value = -2448
